# Standardize filter display - update Debug_Timeline sheet with full asset category breakdown
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Debug_Timeline")

# ---------------------------------------------------------------------------
# 1. Column widths for the new columns D:J (A:C already correctly sized)
# ---------------------------------------------------------------------------
$offset = 5/6
$ws.Columns.Item(4).ColumnWidth  = (20 - $offset)   # D Fondi di investimento
$ws.Columns.Item(5).ColumnWidth  = (13 - $offset)   # E Immobiliare
$ws.Columns.Item(6).ColumnWidth  = (20 - $offset)   # F Liquidità
$ws.Columns.Item(7).ColumnWidth  = (20 - $offset)   # G Oggetti
$ws.Columns.Item(8).ColumnWidth  = (20 - $offset)   # H PAC
$ws.Columns.Item(9).ColumnWidth  = (17 - $offset)   # I Titoli di stato
$ws.Columns.Item(10).ColumnWidth = (20 - $offset)   # J TOTALE

# ---------------------------------------------------------------------------
# 2. Header row: give the new header cells (D1:J1) the same bold style
#    already used by the existing header cells (A1:C1), then set all values.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy($ws.Range("D1:J1"))

$ws.Range("B1").Value = "Criptovalute"
$ws.Range("C1").Value = "ETF"
$ws.Range("D1").Value = "Fondi di investimento"
$ws.Range("E1").Value = "Immobiliare"
$ws.Range("F1").Value = "Liquidità"
$ws.Range("G1").Value = "Oggetti"
$ws.Range("H1").Value = "PAC"
$ws.Range("I1").Value = "Titoli di stato"
$ws.Range("J1").Value = "TOTALE"

$data = New-Object 'object[,]' 29,10
$data[0,0] = "1989-01-16"
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 109180
$data[0,5] = 0
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 0
$data[0,9] = 109180
$data[1,0] = "1997-01-17"
$data[1,1] = 0
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 234318
$data[1,5] = 0
$data[1,6] = 0
$data[1,7] = 0
$data[1,8] = 0
$data[1,9] = 234318
$data[2,0] = "1997-12-01"
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 269592
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 0
$data[2,8] = 0
$data[2,9] = 269592
$data[3,0] = "2004-05-21"
$data[3,1] = 0
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 314592
$data[3,5] = 0
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0
$data[3,9] = 314592
$data[4,0] = "2006-05-22"
$data[4,1] = 0
$data[4,2] = 0
$data[4,3] = 0
$data[4,4] = 341592
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = 0
$data[4,8] = 0
$data[4,9] = 341592
$data[5,0] = "2008-05-22"
$data[5,1] = 0
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 612000
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 0
$data[5,8] = 0
$data[5,9] = 612000
$data[6,0] = "2012-11-30"
$data[6,1] = 0
$data[6,2] = 14391.1
$data[6,3] = 0
$data[6,4] = 612000
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0
$data[6,9] = 626391.1
$data[7,0] = "2014-05-30"
$data[7,1] = 0
$data[7,2] = 17533.30213
$data[7,3] = 0
$data[7,4] = 612000
$data[7,5] = 0
$data[7,6] = 0
$data[7,7] = 0
$data[7,8] = 0
$data[7,9] = 629533.30213
$data[8,0] = "2015-04-08"
$data[8,1] = 0
$data[8,2] = 17533.30213
$data[8,3] = 0
$data[8,4] = 846113
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 0
$data[8,8] = 0
$data[8,9] = 863646.30213
$data[9,0] = "2015-05-22"
$data[9,1] = 0
$data[9,2] = 17533.30213
$data[9,3] = 0
$data[9,4] = 846113
$data[9,5] = 0
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 0
$data[9,9] = 863646.30213
$data[10,0] = "2015-05-27"
$data[10,1] = 0
$data[10,2] = 17533.30213
$data[10,3] = 0
$data[10,4] = 1136488
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 0
$data[10,9] = 1154021.30213
$data[11,0] = "2015-08-04"
$data[11,1] = 0
$data[11,2] = 17533.30213
$data[11,3] = 0
$data[11,4] = 1838827
$data[11,5] = 0
$data[11,6] = 0
$data[11,7] = 0
$data[11,8] = 0
$data[11,9] = 1856360.30213
$data[12,0] = "2019-07-16"
$data[12,1] = 0
$data[12,2] = 17533.30213
$data[12,3] = 6009.006024
$data[12,4] = 1838827
$data[12,5] = 0
$data[12,6] = 0
$data[12,7] = 0
$data[12,8] = 0
$data[12,9] = 1862369.308154
$data[13,0] = "2020-02-11"
$data[13,1] = 0
$data[13,2] = 17533.30213
$data[13,3] = 6009.006024
$data[13,4] = 1838827
$data[13,5] = 0
$data[13,6] = 163008
$data[13,7] = 0
$data[13,8] = 0
$data[13,9] = 2025377.308154
$data[14,0] = "2020-11-02"
$data[14,1] = 0
$data[14,2] = 17533.30213
$data[14,3] = 6009.006024
$data[14,4] = 1838827
$data[14,5] = 0
$data[14,6] = 391000
$data[14,7] = 0
$data[14,8] = 0
$data[14,9] = 2253369.308154
$data[15,0] = "2024-11-13"
$data[15,1] = 0
$data[15,2] = 155263.63205
$data[15,3] = 6009.006024
$data[15,4] = 1838827
$data[15,5] = 0
$data[15,6] = 391000
$data[15,7] = 44998.8521765
$data[15,8] = 0
$data[15,9] = 2436098.4902505
$data[16,0] = "2024-11-14"
$data[16,1] = 0
$data[16,2] = 189097.72205
$data[16,3] = 68987.6346786
$data[16,4] = 1838827
$data[16,5] = 48791.46
$data[16,6] = 391000
$data[16,7] = 44998.8521765
$data[16,8] = 2089.8822
$data[16,9] = 2583792.5511051
$data[17,0] = "2025-03-09"
$data[17,1] = 0
$data[17,2] = 189097.72205
$data[17,3] = 68987.6346786
$data[17,4] = 1846827
$data[17,5] = 48791.46
$data[17,6] = 391000
$data[17,7] = 44998.8521765
$data[17,8] = 2089.8822
$data[17,9] = 2591792.5511051
$data[18,0] = "2025-05-18"
$data[18,1] = 0
$data[18,2] = 200204.85498326
$data[18,3] = 68987.6346786
$data[18,4] = 1846827
$data[18,5] = 48791.46
$data[18,6] = 391000
$data[18,7] = 44998.8521765
$data[18,8] = 2089.8822
$data[18,9] = 2602899.68403836
$data[19,0] = "2025-05-19"
$data[19,1] = 0
$data[19,2] = 200204.85498326
$data[19,3] = 68987.6346786
$data[19,4] = 1846827
$data[19,5] = 48791.46
$data[19,6] = 391000
$data[19,7] = 63861.5517
$data[19,8] = 2089.8822
$data[19,9] = 2621762.38356186
$data[20,0] = "2025-07-31"
$data[20,1] = 0
$data[20,2] = 200204.85498326
$data[20,3] = 68987.6346786
$data[20,4] = 1846827
$data[20,5] = 48791.46
$data[20,6] = 405490.72
$data[20,7] = 63861.5517
$data[20,8] = 2089.8822
$data[20,9] = 2636253.103561861
$data[21,0] = "2025-08-26"
$data[21,1] = 997.027381
$data[21,2] = 209206.216685
$data[21,3] = 68987.6346786
$data[21,4] = 1846827
$data[21,5] = 48791.46
$data[21,6] = 405490.72
$data[21,7] = 63861.5517
$data[21,8] = 2089.8822
$data[21,9] = 2646251.4926446
$data[22,0] = "2025-09-03"
$data[22,1] = 997.027381
$data[22,2] = 209206.216685
$data[22,3] = 68987.6346786
$data[22,4] = 2051827
$data[22,5] = 48791.46
$data[22,6] = 405490.72
$data[22,7] = 63861.5517
$data[22,8] = 2089.8822
$data[22,9] = 2851251.4926446
$data[23,0] = "2025-09-04"
$data[23,1] = 997.027381
$data[23,2] = 215529.714618344
$data[23,3] = 68987.6346786
$data[23,4] = 2270000
$data[23,5] = 48791.46
$data[23,6] = 405490.72
$data[23,7] = 63861.5517
$data[23,8] = 2089.8822
$data[23,9] = 3075747.990577944
$data[24,0] = "2025-09-08"
$data[24,1] = 997.027381
$data[24,2] = 215529.714618344
$data[24,3] = 68987.6346786
$data[24,4] = 2270000
$data[24,5] = 65091.02
$data[24,6] = 405490.72
$data[24,7] = 63861.5517
$data[24,8] = 2089.8822
$data[24,9] = 3092047.550577944
$data[25,0] = "2025-10-01"
$data[25,1] = 997.027381
$data[25,2] = 211014.6991230869
$data[25,3] = 68987.6346786
$data[25,4] = 2270000
$data[25,5] = 65091.02
$data[25,6] = 405490.72
$data[25,7] = 63861.5517
$data[25,8] = 2089.8822
$data[25,9] = 3087532.535082687
$data[26,0] = "2025-10-02"
$data[26,1] = 1073.829225
$data[26,2] = 269674.234384737
$data[26,3] = 84384.43033919999
$data[26,4] = 2270000
$data[26,5] = 65091.02
$data[26,6] = 405490.72
$data[26,7] = 63861.5517
$data[26,8] = 0
$data[26,9] = 3159575.785648937
$data[27,0] = "2025-10-03"
$data[27,1] = 1119.11523402
$data[27,2] = 292283.3118653856
$data[27,3] = 84766.4343524
$data[27,4] = 2270000
$data[27,5] = 70248.14
$data[27,6] = 405490.72
$data[27,7] = 63861.5517
$data[27,8] = 0
$data[27,9] = 3187769.273151806
$data[28,0] = "2025-10-05"
$data[28,1] = 1106.07000798
$data[28,2] = 292607.8108715674
$data[28,3] = 84451.888846638
$data[28,4] = 2270000
$data[28,5] = 70248.14
$data[28,6] = 405490.72
$data[28,7] = 63861.5517
$data[28,8] = 2100.273
$data[28,9] = 3189866.454426186

# ---------------------------------------------------------------------------
# 3. Data rows (A2:J30). Column A holds date-like strings which must stay
#    plain text (not get auto-converted to Excel date serials), so we apply
#    a temporary "Text" number format before assigning values, then clear
#    the formatting afterwards (so the cells end up with the default style,
#    exactly as the rest of the sheet).
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A2:J30")
$ws.Range("A2:A30").NumberFormat = "@"
$dataRange.Value = $data
$dataRange.ClearFormats()

